$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.003.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.49%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.595.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.77%  "
$ws.Range("E3").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'211.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.89%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = "'  -0.10%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.481"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.24%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  +0.52%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  +0.23%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'18.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.07%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.14%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'1.819.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.90%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'1.597.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.91%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("E14").Value = "'  -0.06%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "'  +2.38%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'26.010.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.55%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'60.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.47%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'0.0₃0727"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.64%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  -0.12%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'203.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +5.35%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'  +2.01%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "'  -0.89%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'6.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.84%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +14.42%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'143.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.21%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  -0.07%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "'  -7.24%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  +0.94%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = "'  +1.60%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  +0.82%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  +1.42%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'  +0.32%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'2.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.68%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "'  -0.34%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "'  -0.84%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'1.128.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.88%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  +8.55%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  -0.03%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.796"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.50%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'2.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.07%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.495"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.33%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.776"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.31%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'5.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.36%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'1.733.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.92%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'92.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.08%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "'RenderToken"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.06%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = "'Aave"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'54.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.98%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.0505"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.44%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  +0.58%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'  -0.26%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.0₇0952"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -15.22%  "
$ws.Range("E51").Style = "Normal"
